# Refresh cryptos list prices and 1h volume-change percentages (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.919.39"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "2.919.85"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'594.09"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").Value = "'145.92"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("D9").Value = "'6.80"
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "'0.0000225"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").Value = "'33.68"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "3.401.59"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("D16").Value = "60.911.90"
$ws.Range("E16").Value = "  +0.74%  "
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "2.920.80"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("D19").Value = "'431.06"
$ws.Range("E19").Value = "  +1.75%  "
$ws.Range("D20").Value = "'13.37"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "'0.682"
$ws.Range("E21").Value = "  +1.88%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'81.65"
$ws.Range("E23").Value = "  +1.91%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'11.94"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +4.93%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  -2.98%  "
$ws.Range("D32").Value = "'26.42"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").Value = "0.0₃0850"
$ws.Range("E34").Value = "  +2.32%  "
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("D36").Value = "'5.62"
$ws.Range("E36").Value = "  -0.55%  "
$ws.Range("D37").Value = "'3.03"
$ws.Range("E37").Value = "  +3.18%  "
$ws.Range("D38").Value = "'0.123"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("D40").Value = "'8.57"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").Value = "'0.286"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("E42").Value = "  -2.99%  "
$ws.Range("D43").Value = "'373.40"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "2.702.14"
$ws.Range("E45").Value = "  +2.06%  "
$ws.Range("D46").Value = "'130.98"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D48").Value = "'23.95"
$ws.Range("E48").Value = "  -5.01%  "
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "'0.125"
$ws.Range("E51").Value = "  +2.71%  "
